# Auto-generated edit script applying value updates described by the
# Aegis_Profits.xlsx diff (per-sheet LevePriceNQ/HQ + LeveProfitNQ/HQ
# recompute). Cells whose value changed are set to their new value;
# cells removed by the diff are cleared; cells newly introduced by the
# diff are populated for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 279146.44
$ws.Range("J33").Value = 975863
$ws.Range("L33").Value = 975863
$ws.Range("N33").Value = -976321
$ws.Range("H40").Value = 1979.0555
$ws.Range("J40").Value = 2082.2
$ws.Range("L40").Value = 2082.2
$ws.Range("N40").Value = -2432.2
$ws.Range("H58").Value = 613774.5
$ws.Range("I58").Value = 1032308.75
$ws.Range("J58").Value = 2070.5386
$ws.Range("K58").Value = 3096926.25
$ws.Range("L58").Value = 6211.6158
$ws.Range("M58").Value = -3096776.25
$ws.Range("N58").Value = -6511.6158
$ws.Range("H69").Value = 3464.4443
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("H72").Value = 3464.4443
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("H86").Value = 28808.945
$ws.Range("I86").Value = 43588.332
$ws.Range("J86").Value = 1523.9231
$ws.Range("K86").Value = 43588.332
$ws.Range("L86").Value = 1523.9231
$ws.Range("M86").Value = -42465.332
$ws.Range("N86").Value = -3769.9231
$ws.Range("H89").Value = 28808.945
$ws.Range("I89").Value = 43588.332
$ws.Range("J89").Value = 1523.9231
$ws.Range("K89").Value = 217941.66
$ws.Range("L89").Value = 7619.6155
$ws.Range("M89").Value = -212325.66
$ws.Range("N89").Value = -18851.6155
$ws.Range("H129").Value = 2226.145
$ws.Range("J129").Value = 906.29785
$ws.Range("L129").Value = 2718.89355
$ws.Range("N129").Value = -12718.89355
$ws.Range("H137").Value = 1677.7333
$ws.Range("I137").Value = 1713.8334
$ws.Range("K137").Value = 5141.5002
$ws.Range("M137").Value = -2591.5002
$ws.Range("H138").Value = 1753.9756
$ws.Range("I138").Value = 1851.6
$ws.Range("J138").Value = 1697.6538
$ws.Range("K138").Value = 5554.799999999999
$ws.Range("L138").Value = 5092.9614
$ws.Range("M138").Value = -414.7999999999993
$ws.Range("N138").Value = -15372.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30168.137
$ws.Range("I32").Value = 7610.259
$ws.Range("J32").Value = 193712.75
$ws.Range("K32").Value = 7610.259
$ws.Range("L32").Value = 193712.75
$ws.Range("M32").Value = -7323.259
$ws.Range("N32").Value = -194286.75
$ws.Range("H122").Value = 1165.1562
$ws.Range("I122").Value = 953.8333
$ws.Range("K122").Value = 2861.4999
$ws.Range("M122").Value = -411.4998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808
$ws.Range("H134").Value = 1638.2354
$ws.Range("I134").Value = 1306.7273
$ws.Range("J134").Value = 3722
$ws.Range("K134").Value = 3920.1819
$ws.Range("L134").Value = 11166
$ws.Range("M134").Value = -1385.1819
$ws.Range("N134").Value = -16236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("H31").Value = 33777.688
$ws.Range("I31").Value = 692.75
$ws.Range("J31").Value = 60245.64
$ws.Range("K31").Value = 692.75
$ws.Range("L31").Value = 60245.64
$ws.Range("M31").Value = -397.75
$ws.Range("N31").Value = -60835.64
$ws.Range("H34").Value = 33777.688
$ws.Range("I34").Value = 692.75
$ws.Range("J34").Value = 60245.64
$ws.Range("K34").Value = 692.75
$ws.Range("L34").Value = 60245.64
$ws.Range("M34").Value = -490.75
$ws.Range("N34").Value = -60649.64
$ws.Range("H88").Value = 16947.666
$ws.Range("I88").Value = 5500
$ws.Range("K88").Value = 5500
$ws.Range("M88").Value = -5094
$ws.Range("H91").Value = 16947.666
$ws.Range("I91").Value = 5500
$ws.Range("K91").Value = 5500
$ws.Range("M91").Value = -4096
$ws.Range("H132").Value = 3008.4375
$ws.Range("I132").Value = 2991.5557
$ws.Range("K132").Value = 8974.667099999999
$ws.Range("M132").Value = -6444.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 24000
$ws.Range("J37").Value = 24000
$ws.Range("L37").Value = 72000
$ws.Range("N37").Value = -72224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 12381.546
$ws.Range("J46").Value = 12381.546
$ws.Range("L46").Value = 12381.546
$ws.Range("N46").Value = -12693.546
$ws.Range("H80").Value = 333676670
$ws.Range("I80").Value = 333676670
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 333676670
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -333675672
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 333676670
$ws.Range("I83").Value = 333676670
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 1668383350
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -1668378358
$ws.Range("N83").Value = ""
$ws.Range("H99").Value = 5860.2
$ws.Range("I99").Value = 2955.7778
$ws.Range("J99").Value = 32000
$ws.Range("K99").Value = 2955.7778
$ws.Range("L99").Value = 32000
$ws.Range("M99").Value = -709.7777999999998
$ws.Range("N99").Value = -36492
$ws.Range("H126").Value = 3071
$ws.Range("I126").Value = 3350.375
$ws.Range("J126").Value = 2326
$ws.Range("K126").Value = 10051.125
$ws.Range("L126").Value = 6978
$ws.Range("M126").Value = -7581.125
$ws.Range("N126").Value = -11918
$ws.Range("H132").Value = 3257.8948
$ws.Range("I132").Value = 1749.2
$ws.Range("J132").Value = 4934.222
$ws.Range("K132").Value = 5247.6
$ws.Range("L132").Value = 14802.666
$ws.Range("M132").Value = -2717.6
$ws.Range("N132").Value = -19862.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1069.6
$ws.Range("I22").Value = 1734.2858
$ws.Range("J22").Value = 811.1111
$ws.Range("K22").Value = 1734.2858
$ws.Range("L22").Value = 811.1111
$ws.Range("M22").Value = -1439.2858
$ws.Range("N22").Value = -1401.1111
$ws.Range("H27").Value = 1069.6
$ws.Range("I27").Value = 1734.2858
$ws.Range("J27").Value = 811.1111
$ws.Range("K27").Value = 1734.2858
$ws.Range("L27").Value = 811.1111
$ws.Range("M27").Value = -1627.2858
$ws.Range("N27").Value = -1025.1111
$ws.Range("H82").Value = 1479.2727
$ws.Range("I82").Value = 909
$ws.Range("J82").Value = 2049.5454
$ws.Range("K82").Value = 909
$ws.Range("L82").Value = 2049.5454
$ws.Range("M82").Value = -548
$ws.Range("N82").Value = -2771.5454
$ws.Range("H85").Value = 1479.2727
$ws.Range("I85").Value = 909
$ws.Range("J85").Value = 2049.5454
$ws.Range("K85").Value = 909
$ws.Range("L85").Value = 2049.5454
$ws.Range("M85").Value = 339
$ws.Range("N85").Value = -4545.5454
$ws.Range("H132").Value = 4101.3213
$ws.Range("I132").Value = 4169.1177
$ws.Range("J132").Value = 3996.5454
$ws.Range("K132").Value = 12507.3531
$ws.Range("L132").Value = 11989.6362
$ws.Range("M132").Value = -9977.3531
$ws.Range("N132").Value = -17049.6362
$ws.Range("H136").Value = 1596.3226
$ws.Range("I136").Value = 1391.762
$ws.Range("J136").Value = 2025.9
$ws.Range("K136").Value = 4175.286
$ws.Range("L136").Value = 6077.700000000001
$ws.Range("M136").Value = -1625.286
$ws.Range("N136").Value = -11177.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8781.091
$ws.Range("I15").Value = 7400
$ws.Range("J15").Value = 8919.200000000001
$ws.Range("K15").Value = 7400
$ws.Range("L15").Value = 8919.200000000001
$ws.Range("M15").Value = -7112
$ws.Range("N15").Value = -9495.200000000001
$ws.Range("H62").Value = 9618572
$ws.Range("I62").Value = 38466536
$ws.Range("K62").Value = 38466536
$ws.Range("M62").Value = -38465912
$ws.Range("H65").Value = 9618572
$ws.Range("I65").Value = 38466536
$ws.Range("K65").Value = 192332680
$ws.Range("M65").Value = -192329560
$ws.Range("H81").Value = 333899.5
$ws.Range("I81").Value = 250624.5
$ws.Range("K81").Value = 501249
$ws.Range("M81").Value = -500188
$ws.Range("H84").Value = 333899.5
$ws.Range("I84").Value = 250624.5
$ws.Range("K84").Value = 2506245
$ws.Range("M84").Value = -2500941
$ws.Range("H122").Value = 3285.7273
$ws.Range("I122").Value = 2267
$ws.Range("J122").Value = 3667.75
$ws.Range("K122").Value = 6801
$ws.Range("L122").Value = 11003.25
$ws.Range("M122").Value = -4351
$ws.Range("N122").Value = -15903.25
$ws.Range("H132").Value = 2184.5615
$ws.Range("I132").Value = 2247.6
$ws.Range("J132").Value = 2036.2354
$ws.Range("K132").Value = 6742.799999999999
$ws.Range("L132").Value = 6108.706200000001
$ws.Range("M132").Value = -4212.799999999999
$ws.Range("N132").Value = -11168.7062
$ws.Range("H136").Value = 1492.4286
$ws.Range("I136").Value = 592.9643
$ws.Range("J136").Value = 2212
$ws.Range("K136").Value = 1778.8929
$ws.Range("L136").Value = 6636
$ws.Range("M136").Value = 771.1071000000002
$ws.Range("N136").Value = -11736
